$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "展览" (Exhibitions) - bump "want to go" counts (column F)
# ---------------------------------------------------------------------------
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F4").Value  = 1262
$wsExpo.Range("F5").Value  = 61
$wsExpo.Range("F7").Value  = 960
$wsExpo.Range("F8").Value  = 927
$wsExpo.Range("F14").Value = 1787
$wsExpo.Range("F15").Value = 3741
$wsExpo.Range("F16").Value = 1121
$wsExpo.Range("F18").Value = 2531
$wsExpo.Range("F20").Value = 1058
$wsExpo.Range("F21").Value = 3488
$wsExpo.Range("F23").Value = 828
$wsExpo.Range("F24").Value = 33
$wsExpo.Range("F25").Value = 2084
$wsExpo.Range("F26").Value = 105
$wsExpo.Range("F27").Value = 809
$wsExpo.Range("F29").Value = 161
$wsExpo.Range("F30").Value = 174
$wsExpo.Range("F32").Value = 1290
$wsExpo.Range("F33").Value = 1910
$wsExpo.Range("F34").Value = 472
$wsExpo.Range("F36").Value = 579
$wsExpo.Range("F37").Value = 259
$wsExpo.Range("F40").Value = 72

# ---------------------------------------------------------------------------
# Sheet "演出" (Performances)
# ---------------------------------------------------------------------------
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F10").Value = 8
$wsShow.Range("G16").Value = 1280

# ---------------------------------------------------------------------------
# Sheet "本地生活" (Local life)
# ---------------------------------------------------------------------------
$wsLocal = $wb.Worksheets.Item("本地生活")
$wsLocal.Range("F2").Value = 293

# ---------------------------------------------------------------------------
# Sheet "全部类型" (All types) - mirrors the other sheets plus its own rows
# ---------------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value  = 1262
$wsAll.Range("F4").Value  = 61
$wsAll.Range("F5").Value  = 960
$wsAll.Range("F6").Value  = 927
$wsAll.Range("F15").Value = 1787
$wsAll.Range("F16").Value = 3741
$wsAll.Range("F17").Value = 1121
$wsAll.Range("F20").Value = 2531
$wsAll.Range("F22").Value = 1058
$wsAll.Range("F23").Value = 3488
$wsAll.Range("F25").Value = 828
$wsAll.Range("F27").Value = 33
$wsAll.Range("F28").Value = 2084
$wsAll.Range("F30").Value = 8
$wsAll.Range("F32").Value = 105
$wsAll.Range("F34").Value = 809
$wsAll.Range("F36").Value = 161
$wsAll.Range("F37").Value = 174
$wsAll.Range("F40").Value = 1290
$wsAll.Range("F41").Value = 1910

# Row 42 on "全部类型" is replaced with a different event entirely
# (leading apostrophe forces text so Excel doesn't auto-convert the
# "yyyy-mm-dd" looking string into a date serial; Style reset afterwards
# drops the quote-prefix formatting Excel applies so the cell keeps its
# original, unstyled look - matching the source inlineStr cell.)
$wsAll.Range("B42").Value = "'2024-06-14"
$wsAll.Range("B42").Style = "Normal"
$wsAll.Range("C42").Value = "杭州·苗阜王声 青曲社相声全国巡演"
$wsAll.Range("D42").Value = "湖墅南路138号 杭州浙话艺术剧院"
$wsAll.Range("E42").Value = "2024.06.14 19:30-06.14 22:00"
$wsAll.Range("F42").Value = 2
$wsAll.Range("G42").Value = 280
$wsAll.Range("H42").Value = "https://show.bilibili.com/platform/detail.html?id=83382"
$wsAll.Range("I42").Value = "//i1.hdslb.com/bfs/openplatform/202403/hUGL3xz01711346789039.jpeg"

$wsAll.Range("G43").Value = 1280
$wsAll.Range("F44").Value = 472
$wsAll.Range("F45").Value = 579
$wsAll.Range("F46").Value = 259
$wsAll.Range("F49").Value = 72
